# OW-268 link uploaded trades to their account
#
# The "Position Account ID" column (B) on the IRS-Cleared sheet held a
# placeholder FCM/clearing-member account number ("MEGA104"). Update it to
# the uploaded trade's linked account ("acc1"), then leave the selection on
# the edited cell (mirrors what happens when a user types the new value into
# B2 in the Excel UI and the workbook is saved with that as the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "acc1"
$ws.Range("B2").Select()
